# Update "想去人数" (want-to-go count) values in column F across the four
# sheets of the workbook, matching the regenerated data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1086
$ws1.Range("F3").Value = 4652
$ws1.Range("F4").Value = 601
$ws1.Range("F6").Value = 1822
$ws1.Range("F7").Value = 43
$ws1.Range("F8").Value = 736
$ws1.Range("F9").Value = 35
$ws1.Range("F12").Value = 1129
$ws1.Range("F14").Value = 808
$ws1.Range("F15").Value = 1015
$ws1.Range("F16").Value = 550
$ws1.Range("F17").Value = 516
$ws1.Range("F18").Value = 624
$ws1.Range("F19").Value = 171
$ws1.Range("F21").Value = 1187
$ws1.Range("F23").Value = 2509
$ws1.Range("F25").Value = 1556
$ws1.Range("F26").Value = 488
$ws1.Range("F29").Value = 4246

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 15
$ws2.Range("F7").Value = 400
$ws2.Range("F9").Value = 4158
$ws2.Range("F11").Value = 26
$ws2.Range("F14").Value = 16
$ws2.Range("F35").Value = 32

# Sheet "本地生活"
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 1731
$ws3.Range("F6").Value = 1080
$ws3.Range("F7").Value = 247

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1731
$ws4.Range("F4").Value = 1080
$ws4.Range("F5").Value = 247
$ws4.Range("F6").Value = 15
$ws4.Range("F7").Value = 1086
$ws4.Range("F9").Value = 4652
$ws4.Range("F10").Value = 601
$ws4.Range("F12").Value = 1822
$ws4.Range("F13").Value = 736
$ws4.Range("F18").Value = 1129
$ws4.Range("F20").Value = 26
$ws4.Range("F22").Value = 808
$ws4.Range("F23").Value = 1015
$ws4.Range("F24").Value = 550
$ws4.Range("F25").Value = 516
$ws4.Range("F26").Value = 624
$ws4.Range("F27").Value = 171
$ws4.Range("F33").Value = 1187
$ws4.Range("F37").Value = 2509
$ws4.Range("F43").Value = 1556
$ws4.Range("F44").Value = 488
$ws4.Range("F48").Value = 4246
$ws4.Range("F49").Value = 32
